# Apply "add both yard blocks" edit to GreenLine_Layout.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: Row 152 - new yard block for section 57 (block 151) ---
# Apply the same "centered" look used by the rest of the table (style 4)
# to the normal data columns, and the centered 0.00-format look (style 5)
# to the Time-To-Traverse column, before filling in the values.
$ws1.Range("B152:H152").HorizontalAlignment = -4108
$ws1.Range("M152:O152").HorizontalAlignment = -4108
$ws1.Range("Q152").HorizontalAlignment = -4108
$ws1.Range("S152").HorizontalAlignment = -4108
$ws1.Range("R152").HorizontalAlignment = -4108
$ws1.Range("R152").NumberFormat = "0.00"

$ws1.Range("B152").Value = "$"
$ws1.Range("C152").Value = 151
$ws1.Range("D152").Value = 100
$ws1.Range("E152").Value = 0
$ws1.Range("F152").Value = 15
$ws1.Range("G152").Value = 0
$ws1.Range("H152").Value = 0.5
$ws1.Range("L152").Value = "K63"
$ws1.Range("M152").Value = 1
$ws1.Range("N152").Value = 0
$ws1.Range("O152").Value = 0
$ws1.Range("Q152").Value = 0
$ws1.Range("R152").Value = 24
$ws1.Range("S152").Value = 2
$ws1.Range("T152").Value = 0
$ws1.Range("U152").Formula = "=F152*0.621371"
$ws1.Range("V152").Formula = "=D152*1.09361"

# --- Sheet1: Row 153 - new yard block for section 63 (block 152) ---
$ws1.Range("B153:H153").HorizontalAlignment = -4108
$ws1.Range("M153:O153").HorizontalAlignment = -4108
$ws1.Range("Q153").HorizontalAlignment = -4108
$ws1.Range("S153").HorizontalAlignment = -4108
$ws1.Range("R153").HorizontalAlignment = -4108
$ws1.Range("R153").NumberFormat = "0.00"

$ws1.Range("B153").Value = "y"
$ws1.Range("C153").Value = 152
$ws1.Range("D153").Value = 100
$ws1.Range("E153").Value = 0
$ws1.Range("F153").Value = 15
$ws1.Range("G153").Value = 0
$ws1.Range("H153").Value = 0.5
$ws1.Range("L153").Value = "I57"
$ws1.Range("M153").Value = 1
$ws1.Range("N153").Value = 0
$ws1.Range("O153").Value = 0
$ws1.Range("Q153").Value = 0
$ws1.Range("R153").Value = 24
$ws1.Range("S153").Value = 2
$ws1.Range("T153").Value = 0
$ws1.Range("U153").Formula = "=F153*0.621371"
$ws1.Range("V153").Formula = "=D153*1.09361"

# --- Sheet1: update the two "Yard" exit-block labels in column K ---
# K58 referenced block 57's yard exit: "57-Yard, 57-58" -> "57-152, 57-58"
# K64 referenced block 63's yard exit: "63-Yard, 63-62" -> "63-151, 63-62"
$ws1.Range("K58").Value = "57-152, 57-58"
$ws1.Range("K64").Value = "63-151, 63-62"

# --- Sheet2: add per-section block counts for the two new sections ---
$ws2.Range("A28").Value = "$"
$ws2.Range("B28").Value = 1
$ws2.Range("A29").Value = "y"
$ws2.Range("B29").Value = 1

# --- Sheet2: update the view state to match the edited file ---
$ws2.Range("B33").Select() | Out-Null

# --- Sheet1: re-activate sheet1 (it is the tab shown on open) and
#     update its view state to match the edited file ---
$ws1.Activate()
$ws1.Range("K156").Select() | Out-Null
